$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.546.48"
$ws.Range("E2").Value = "  -0.23%  "
$ws.Range("D3").Value = "1.812.47"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.69"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("E6").Value = "  +4.09%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  +6.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.303"
$ws.Range("E9").Value = "  +2.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0698"
$ws.Range("E10").Value = "  +0.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0958"
$ws.Range("E11").Value = "  +0.88%  "
$ws.Range("D12").Value = "2.073.45"
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.25"
$ws.Range("E13").Value = "  +0.62%  "
$ws.Range("D14").Value = "1.826.45"
$ws.Range("E14").Value = "  +1.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.651"
$ws.Range("E15").Value = "  +1.81%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.52"
$ws.Range("E16").Value = "  +4.32%  "
$ws.Range("D17").Value = "34.503.12"
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.29"
$ws.Range("E18").Value = "  +0.45%  "
$ws.Range("D19").Value = "0.0₃0800"
$ws.Range("E19").Value = "  -0.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "246.38"
$ws.Range("E20").Value = "  -0.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.46"
$ws.Range("E21").Value = "  +1.04%  "
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "172.33"
$ws.Range("E24").Value = "  +0.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.12"
$ws.Range("E25").Value = "  +1.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.08"
$ws.Range("E26").Value = "  +10.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.85"
$ws.Range("E27").Value = "  +1.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.120"
$ws.Range("E28").Value = "  +2.83%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.05"
$ws.Range("E30").Value = "  -0.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0536"
$ws.Range("E31").Value = "  +1.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.87"
$ws.Range("E32").Value = "  +1.52%  "
$ws.Range("E33").Value = "  +0.63%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.85"
$ws.Range("E34").Value = "  +0.38%  "
$ws.Range("D35").Value = "1.397.66"
$ws.Range("E35").Value = "  -2.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.683"
$ws.Range("E36").Value = "  +0.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.47"
$ws.Range("E37").Value = "  -3.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.08"
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "83.95"
$ws.Range("E40").Value = "  -1.14%  "
$ws.Range("E41").Value = "  +1.67%  "
$ws.Range("E42").Value = "  +2.80%  "
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.31"
$ws.Range("E44").Value = "  -3.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.13"
$ws.Range("E45").Value = "  +4.59%  "
$ws.Range("E46").Value = "  -2.09%  "
$ws.Range("E47").Value = "  -1.88%  "
$ws.Range("D48").Value = "1.973.64"
$ws.Range("E48").Value = "  +0.50%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.37"
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("E50").Value = "  +2.02%  "
$ws.Range("E51").Value = "  +0.14%  "
